# "Generate Report for Handback" - records that a handback was produced for
# 035cb385-fb23-425e-9539-1a938b2a33c5.md in both the zh-cn and de-de sheets,
# but flags that the handback file version is stale (not the latest commit).

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bccc99c9dd630c62a7e8508795f8c88d0c5e2733/e2e/035cb385-fb23-425e-9539-1a938b2a33c5.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6c10a184334c95831a813c5fd4eefd9020278aad/e2e/035cb385-fb23-425e-9539-1a938b2a33c5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bccc99c9dd630c62a7e8508795f8c88d0c5e2733/e2e/035cb385-fb23-425e-9539-1a938b2a33c5.md."
$targetFileValue = "035cb385-fb23-425e-9539-1a938b2a33c5.md"
$hyperlinkFontColor = 15570276

function Apply-HandbackRow($SheetName, $HandbackFileValue, $HandbackDateTime) {
    $ws = $wb.Worksheets.Item($SheetName)

    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(16).ColumnWidth = 39.166666666666664

    $ws.Range("I5").Value = $targetFileValue
    $ws.Range("I5").Font.Color = $hyperlinkFontColor
    $ws.Range("I5").Font.Underline = 2
    $ws.Hyperlinks.Add($ws.Range("I5"), $targetUrl, "", "", $targetFileValue) | Out-Null

    $ws.Range("J5").Value = $HandbackFileValue
    $ws.Range("K5").Value = $HandbackDateTime
    $ws.Range("P5").Value = $errorDetail
}

Apply-HandbackRow "zh-cn" "035cb385-fb23-425e-9539-1a938b2a33c5.ad23e57dc60e4b4fea4560df0d31a43c28572975.zh-cn.xlf" "2016-10-18 12:07:02"
Apply-HandbackRow "de-de" "035cb385-fb23-425e-9539-1a938b2a33c5.ad23e57dc60e4b4fea4560df0d31a43c28572975.de-de.xlf" "2016-10-18 12:07:19"
